# Update VoiceImage to StreetImage.
#
# Applies the Timelist.xlsx changes:
#  - Append wording to the last changelog entry (shared string).
#  - Change B23 formula from 1+1 to 1+1+1 (and resulting SUM in B27).
#  - Add three new date rows (24-26) following the same style as row 23.
#  - Move the sheet's active selection from C24 to B24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the changelog text for the last entry (C23 / shared string).
$ws.Range("C23").Value = "Reduce GPU computation when loading GIF. Create VoicePFImage to handle image loading. AssistantHorizontalView works now."

# 2) Update the B23 formula (1+1 -> 1+1+1); B27's SUM recalculates automatically.
$ws.Range("B23").Formula = "=1+1+1"

# 3) Add three new rows continuing the date sequence after row 23.
$ws.Range("A24").Value = 41951
$ws.Range("A25").Value = 41952
$ws.Range("A26").Value = 41953

# Copy the formatting/style from A23 onto the new date cells so they match
# the existing date column style (s="3") instead of the default column style.
$ws.Range("A23").Copy()
[void]$ws.Range("A24:A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) Move the active selection from C24 to B24.
[void]$ws.Range("B24").Select()
